$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - account holder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number and surname
# Card number looks numeric, but must stay text (same style as before) -
# enter with a leading apostrophe then re-apply the original cell format
# only (no value change) so Excel doesn't stick a quote-prefix style on it.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 26.01.2024"

# Row 6
$ws.Range("B6").Value = "28.01."
$ws.Range("C6").Value = "29.01."
$ws.Range("E6").Value = "25,07-"

# Row 7
$ws.Range("B7").Value = "31.01."
$ws.Range("C7").Value = "01.02."
$ws.Range("D7").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E7").Value = "88,52-"

# Row 8
$ws.Range("B8").Value = "01.02."
$ws.Range("C8").Value = "02.02."
$ws.Range("D8").Value = "PAYPAL KKVOLX"
$ws.Range("E8").Value = "5,65-"

# Row 9
$ws.Range("B9").Value = "04.02."
$ws.Range("C9").Value = "05.02."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-17969281"
$ws.Range("E9").Value = "55,21-"

# Row 10
$ws.Range("B10").Value = "05.02."
$ws.Range("C10").Value = "06.02."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 90406919"
$ws.Range("E10").Value = "40,41-"

# Row 11 - the transaction is removed entirely, leaving the row blank.
# E11 picks up the wrap/right/vertical-center alignment used by the blank
# "spacer" rows elsewhere in this template.
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").WrapText = $true
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").HorizontalAlignment = -4152

# Row 12 - closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 10.02.2024"
$ws.Range("E12").Value = "214,86-"

# Row 13 - next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.02.2024"
